$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new "sala" rows (10-12) with name, branch, capacity-ish numbers
$ws.Range("A10").Value = "SALAMax1"
$ws.Range("B10").Value = "Sucursal Cartago"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 80

$ws.Range("A11").Value = "SALA123"
$ws.Range("B11").Value = "Sucursal Cartago"
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 80

$ws.Range("A12").Value = "sala1"
$ws.Range("B12").Value = "Sucursal Cartago"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 80

# Touch page setup so the worksheet emits a <headerFooter/> element, as in
# the authored edit (Excel writes this once header/footer settings are
# touched through the UI/object model).
$ws.PageSetup.CenterHeader = ""
